# [TASK] Working on documentation
# Add newly logged time-tracking entries to the "Stundenerfassung" sheet
# and correct a typo in an existing entry's description.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stundenerfassung")

# --- Row 4: 2019-07-16, MG, Morphologische Analyse, 1h ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null
$ws.Range("A4").Value = 43662
$ws.Range("B4").Value = "MG"
$ws.Range("C4").Value = "Morphologische Analyse"
$ws.Range("D4").Value = 1

# --- Row 5: 2019-07-19, MG, Morphologische Analyse, 2h ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null
$ws.Range("A5").Value = 43665
$ws.Range("B5").Value = "MG"
$ws.Range("C5").Value = "Morphologische Analyse"
$ws.Range("D5").Value = 2

# --- Row 6: 2019-07-19, MG, Risk Managment, 2.5h ---
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null
$ws.Range("A6").Value = 43665
$ws.Range("B6").Value = "MG"
$ws.Range("C6").Value = "Risk Managment"
$ws.Range("D6").Value = 2.5

$excel.CutCopyMode = $false

# --- Fix typo in the first logged entry's description ---
# (done last so new shared-string entries are appended in the same
#  order Excel would have produced them)
$ws.Range("C2").Value = "Aufsetzen Projekt auf Github"

$ws.Range("C2").Select() | Out-Null
